$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.575
$ws.Range("C4").Value = -12.566
$ws.Range("A6").Value = -22.176
$ws.Range("A7").Value = -19.81
$ws.Range("A8").Value = -21.961
$ws.Range("C8").Value = -12.404
$ws.Range("C9").Value = -11.361
$ws.Range("C12").Value = -10.919
$ws.Range("A16").Value = -22.037
$ws.Range("C17").Value = -13.297
$ws.Range("C18").Value = -11.591
$ws.Range("C19").Value = -12.126
$ws.Range("A20").Value = -19.828
$ws.Range("C20").Value = -11.729
$ws.Range("A21").Value = -20.269
$ws.Range("C26").Value = -12.046
$ws.Range("A28").Value = -22.021
$ws.Range("A29").Value = -21.344
$ws.Range("A30").Value = -21.917
$ws.Range("C31").Value = -13.298
$ws.Range("A32").Value = -21.832
$ws.Range("C39").Value = -12.117
$ws.Range("A40").Value = -19.922
$ws.Range("C40").Value = -12.226
$ws.Range("C41").Value = -12.117
$ws.Range("C42").Value = -12.14
$ws.Range("C43").Value = -12.128
$ws.Range("A46").Value = -21.89
$ws.Range("C47").Value = -11.998
$ws.Range("C48").Value = -11.74
$ws.Range("A51").Value = -21.715
$ws.Range("A52").Value = -21.972
$ws.Range("C54").Value = -13.268
$ws.Range("A57").Value = -22.275
$ws.Range("A59").Value = -22.429
$ws.Range("A62").Value = -21.97
$ws.Range("C62").Value = -13.419
$ws.Range("C63").Value = -11.178
$ws.Range("C64").Value = -11.011
$ws.Range("A66").Value = -21.567
$ws.Range("A73").Value = -20.387
$ws.Range("A74").Value = -21.036
$ws.Range("C76").Value = -12.805
$ws.Range("A77").Value = -20.407
$ws.Range("C81").Value = -12.853
$ws.Range("C84").Value = -13.616
$ws.Range("C89").Value = -13.612
$ws.Range("A92").Value = -21.759
$ws.Range("C94").Value = -11.539
$ws.Range("A100").Value = -22.195